$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Row 10 (Objetivos:): replace the (wrong) B/C text with the real objectives paragraph.
$objetivos = "Experiências em laboratório de caráter multidisciplinar que tem por objetivo colocar o aluno em contato com equipamentos de engenharia e consolidar os conceitos de fenômenos de transporte. O desenvolvimento das atividades inclui montagem, medidas e interpretação de resultados em áreas relevantes da engenharia como cinética e reatores químicos, fenômenos de transporte, operações unitárias e processos químicos industriais. A disciplina permite um programa dinâmico, onde os experimentos poderão ser mudados e/ou revezados em função da evolução dos laboratórios ou necessidades específicas."
$ws.Range("B10").Value = $objetivos
$ws.Range("C10").Value = $objetivos

# 2) Insert a new row at 13 (pushes the old rows 13-24 down to 14-25, carrying
#    their content/heights along automatically).
$ws.Rows("13:13").Insert()

# The inserted row 13 picked up column-A's (bold, non-wrap) formatting; give
# B13/C13 the correct column formatting by copying it from row 14 (which still
# has the original formats for columns B/C), then set the values and drop the
# stray A13 formatting.
$ws.Range("B14").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$docentes = "5816812 - João Paulo Alves Silva"
$ws.Range("B13").Value = $docentes
$ws.Range("C13").Value = $docentes
$ws.Range("A13").Clear()

# 3) Update the text for the rows that now hold new/different content
#    (labels stayed put; only the B/C bodies change).
$resumido = "1) Reatores químicos`n2) Fermentação`n3) Processos químicos"
$ws.Range("B14").Value = $resumido
$ws.Range("C14").Value = $resumido

$programa = "1) Reatores químicos: operação de reator de mistura para avaliar a influência do tempo de residência na conversão.`n2) Fermentação: determinação de parâmetros cinéticos em processo de fermentação alcoólica por leveduras em reatores bioquímicos.`n3) Processos químicos: executar síntese de produto de interesse industrial em uma planta de pequeno porte. Poderão ser observados aspectos como instrumentação (controles de vazão, temperatura, nível, etc), equipamentos diversos de operações unitárias, sistema de aquisição de dados, etc"
$ws.Range("B16").Value = $programa
$ws.Range("C16").Value = $programa

$metodo = "Aplicação de prova(s) e relatório(s)."
$ws.Range("B19").Value = $metodo
$ws.Range("C19").Value = $metodo

$criterio = "A média do período será definida pelo professor da disciplina. Alunos com média final igual ou superior a 5,0 estarão aprovados, desde que tenham freqüência mínima de 70% (regimental). Alunos com média inferior a 3,0 e/ou freqüência inferior a 70% estarão reprovados (regimental). Alunos com média superior ou igual a 3,0 e inferior a 5,0 e que tenham freqüência mínima de 70% serão submetidos ao período de recuperação (regimental)."
$ws.Range("B20").Value = $criterio
$ws.Range("C20").Value = $criterio

$norma = "A média final após a recuperação para a disciplina será a média aritmética entre a média do período e a nota da recuperação. Durante o período de recuperação, poderá ser marcada uma aula com a finalidade de sanar dúvidas e/ou revisar conceitos fundamentais. Em data posterior os alunos serão submetidos a uma prova de recuperação"
$ws.Range("B21").Value = $norma
$ws.Range("C21").Value = $norma

$biblio = "1) FOUST, Alan S.; WENZEL, Leonard A.; CLUMP, Curtis W.; MAUS, Louis; ANDERSEN, L. Bryce. Princípios das operações unitárias. Rio de Janeiro: Guanabara Dois/LTC, 1982.`n2) GEANKOPLIS, Christie John. Transport Processes and Separation Process Principles. New York: Prentice Hall, 2003.`n3) COUPER, James R.; PENNEY, W. Roy; FAIR, James R.; WALAS, Stanley M. Chemical Process Equipment: Selection and Design. Amsterdam: Elsevier, 2005.`n4) FOGLER, H. S. Elementos de engenharia das reações químicas. 3.ed. Rio de Janeiro: LTC Editora, 2002.`n5) LEVENSPIEL, O. Chemical Reaction Engineering. 3rd.ed. New York: John Wiley & Sons, 1998.`n6) PERRY, Robert H.; GREEN, Don W. Perry's Chemical Engineers' Handbook. 8th.ed. New York: McGraw-Hill, 2008."
$ws.Range("B22").Value = $biblio
$ws.Range("C22").Value = $biblio
